$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.464.73'
$ws.Range('E2').Value = '  +1.74%  '
$ws.Range('D3').Value = '2.659.94'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '608.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.31'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.08%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.02%  '
$ws.Range('E9').Value = '  +1.56%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.387'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +6.48%  '
$ws.Range('E11').Value = '  -0.67%  '
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('E13').Value = '  +3.00%  '
$ws.Range('D14').Value = '3.136.89'
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').Value = '64.291.38'
$ws.Range('E15').Value = '  +1.71%  '
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').Value = '2.655.59'
$ws.Range('E17').Value = '  +0.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.17'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.72%  '
$ws.Range('E19').Value = '  +4.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '350.93'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.39%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.45%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.81'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('E25').Value = '  +13.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.41'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.48%  '
$ws.Range('E27').Value = '  +4.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.21'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '555.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.164'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('E32').Value = '  +1.44%  '
$ws.Range('D33').Value = '0.0₃0867'
$ws.Range('E33').Value = '  +6.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.77'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.37'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '168.56'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.22%  '
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('E38').Value = '  +7.76%  '
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.45'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.57%  '
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '166.98'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.25'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.86'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0578'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '22.10'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.50%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('E48').Value = '  +15.31%  '
$ws.Range('E49').Value = '  +3.18%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0968'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.12'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.45%  '
